# Generate Report for Handoff
#
# The handback for "5ad28dad-3920-4125-bac1-dd1defa03890.md" is no longer
# the latest version, so its status flips from "Handed back: in sync with
# en-US" to "Ready for handoff" everywhere it is reported (the Overview
# rollup sheet plus each per-locale detail sheet), the relevant
# timestamps are refreshed, and an explanatory error message is recorded
# for each locale.

$wb = $excel.ActiveWorkbook

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a5183b40492ef0ceafda257f1fb4c75a6ce7c300/e2e/5ad28dad-3920-4125-bac1-dd1defa03890.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/19abc2ec4a5065af78be0a1d77548baa38fba0df/e2e/5ad28dad-3920-4125-bac1-dd1defa03890.md."

# --- Overview sheet: row 3 is the 5ad28dad-... file -----------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-08-31 15:02:23"

# --- zh-cn detail sheet: row 3 is the 5ad28dad-... file --------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("H3").Value = "2016-08-31 15:02:10"
$wsZhCn.Range("P3").Value = $errorDetail
$wsZhCn.Range("P1").ColumnWidth = 39.17

# --- de-de detail sheet: row 3 is the 5ad28dad-... file --------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("H3").Value = "2016-08-31 15:02:23"
$wsDeDe.Range("P3").Value = $errorDetail
$wsDeDe.Range("P1").ColumnWidth = 39.17
